$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header labels in B1:H1, keeping A1's "name" label intact.
# (ClearContents() drops the value but keeps each cell's existing style.)
$ws.Range("B1:H1").ClearContents()

# E1, F1 and H1 were the wrapped header cells - re-assert WrapText explicitly
# so the formatting survives even though the text driving it is now gone.
$ws.Range("E1:F1").WrapText = $true
$ws.Range("H1").WrapText = $true

# A1 ("name") becomes left-aligned with the font explicitly applied.
$ws.Range("A1").HorizontalAlignment = -4131  # xlLeft

# Move the active selection to C1.
$ws.Range("C1").Select() | Out-Null
